$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table was refreshed with corrected department names/values
# (commit message: "cambio de nombre a files" - "los files tenian nombres incorrectos")
# Columns: A=Pais, B=Departamento, C=Muestras, D=Confirmado (+), E=% de Positividad

$rows = @(
    @{ Row = 2;  Dept = "LIMA";           C = 224601; D = 32339; E = 14.399999999999999 },
    @{ Row = 3;  Dept = "CALLAO";         C = 18164;  D = 4245;  E = 23.369999999999997 },
    @{ Row = 4;  Dept = "AREQUIPA";       C = 15613;  D = 763;   E = 4.8899999999999997 },
    @{ Row = 5;  Dept = "LAMBAYEQUE";     C = 14165;  D = 3008;  E = 21.240000000000002 },
    @{ Row = 6;  Dept = "PIURA";          C = 13599;  D = 1804;  E = 13.270000000000001 },
    @{ Row = 7;  Dept = "LA LIBERTAD";    C = 12503;  D = 1062;  E = 8.49 },
    @{ Row = 8;  Dept = "ANCASH";         C = 12098;  D = 1159;  E = 9.58 },
    @{ Row = 9;  Dept = "JUNIN";          C = 8389;   D = 587;   E = 7.0000000000000009 },
    @{ Row = 10; Dept = "LORETO";         C = 7784;   D = 1595;  E = 20.49 },
    @{ Row = 11; Dept = "ICA";            C = 7300;   D = 698;   E = 9.56 },
    @{ Row = 12; Dept = "SAN MARTIN";     C = 6487;   D = 303;   E = 4.67 },
    @{ Row = 13; Dept = "CAJAMARCA";      C = 6271;   D = 322;   E = 5.13 },
    @{ Row = 14; Dept = "CUSCO";          C = 6162;   D = 249;   E = 4.04 },
    @{ Row = 15; Dept = "TACNA";          C = 6153;   D = 149;   E = 2.42 },
    @{ Row = 16; Dept = "PUNO";           C = 5952;   D = 144;   E = 2.42 },
    @{ Row = 17; Dept = "MADRE DE DIOS";  C = 5652;   D = 124;   E = 2.19 },
    @{ Row = 18; Dept = "UCAYALI";        C = 5271;   D = 1032;  E = 19.580000000000002 },
    @{ Row = 19; Dept = "HUANUCO";        C = 4782;   D = 255;   E = 5.33 },
    @{ Row = 20; Dept = "MOQUEGUA";       C = 4421;   D = 142;   E = 3.2099999999999995 },
    @{ Row = 21; Dept = "TUMBES";         C = 3927;   D = 378;   E = 9.629999999999999 },
    @{ Row = 22; Dept = "HUANCAVELICA";   C = 3701;   D = 203;   E = 5.4899999999999993 },
    @{ Row = 23; Dept = "APURIMAC";       C = 3694;   D = 111;   E = 3 },
    @{ Row = 24; Dept = "AYACUCHO";       C = 3576;   D = 181;   E = 5.0599999999999996 },
    @{ Row = 25; Dept = "AMAZONAS";       C = 3203;   D = 170;   E = 5.3100000000000005 },
    @{ Row = 26; Dept = "PASCO";          C = 3111;   D = 166;   E = 5.34 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value2 = $r.Dept
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D
    $ws.Cells.Item($r.Row, 5).Value2 = $r.E
}

# The % column is no longer formatted as a percentage (plain numbers now)
$ws.Range("E2:E26").Style = "Normal"

# Update the selected cell shown in the saved view
[void]$ws.Range("J10").Select()
